$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-5 get their D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg)
# values cyclically rotated: row2->row3, row3->row5, row5->row4, row4->row2.
# Capture the original values first so the writes don't clobber each other.

$rows = @(2, 3, 4, 5)
$cols = @("D", "M", "N", "O", "P", "S")

$original = @{}
foreach ($r in $rows) {
    $original[$r] = @{}
    foreach ($c in $cols) {
        $original[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row -> source row mapping (i.e. new row gets old value of source row)
$mapping = @{ 2 = 4; 3 = 2; 4 = 5; 5 = 3 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $original[$src][$c]
    }
}
